# Add 2022-Q1 data
# ------------------------------------------------------------------
# 1) Insert a new sheet "2022-Q1" right before the "总计" (totals) sheet,
#    cloning the layout/format used by the existing quarterly sheets
#    (e.g. "2021-Q4") since they all share the same column layout.
# 2) Populate it with the single fund holding row for 2022-Q1.
# 3) Insert a new leading row into "总计" for the 2022-Q1 aggregate and
#    renumber the existing index column.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")

# --- 1) New sheet, positioned right before "总计" -------------------
# NOTE: passing a worksheet object as the Add() "Before" arg repoints
# *that same variable handle* at the freshly created sheet (an
# iron_native COM quirk) -- so resolve "总计" fresh (by name) any time
# it's needed instead of caching it in a variable across the Add() call.
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Copy the header/body formatting (fonts, borders, alignment) from the
# 2021-Q4 sheet so the new sheet matches the existing look exactly.
# (Deliberately excludes A1, which has no cell/formatting in the source
# sheets -- copying it would create a spurious empty <c r="A1"/>.)
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:H2").Copy()
$newSheet.Range("A2:H2").PasteSpecial(-4122)

# --- 2) Header row ---------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- Data row (text-typed values, matching source formatting) -------
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'510810"
$newSheet.Range("C2").Value = "汇添富中证上海国企ETF"
$newSheet.Range("D2").Value = "'68.43"
$newSheet.Range("E2").Value = "'99.71"
$newSheet.Range("F2").Value = "'3.16"
$newSheet.Range("G2").Value = "'2.1624"
$newSheet.Range("H2").Value = 9

# --- 3) Update "总计" sheet: insert 2022-Q1 as the new first data row
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Re-apply the standard row formatting (lost/odd after the row insert)
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 2.16

# Renumber the index column (A) for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# Values for the shifted rows stay the same (2021-Q4, 2021-Q3, ... 2020-Q4),
# only their row position changed during the insert, so no other edits
# are required.
